$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-50: ticker symbol changes ---
$ws.Range("B2").Value = "NSE:5PAISA"
$ws.Range("C2").Value = "NSE:BANG"
$ws.Range("D2").Value = "NSE:EXIDEIND"
$ws.Range("F2").Value = "NSE:ABCAPITAL"
$ws.Range("B3").Value = "NSE:AAKASH"
$ws.Range("C3").Value = "NSE:HPL"
$ws.Range("D3").Value = "NSE:INDUSTOWER"
$ws.Range("F3").Value = "NSE:JIOFIN"
$ws.Range("B4").Value = "NSE:ABCAPITAL"
$ws.Range("C4").Value = "NSE:KSCL"
$ws.Range("D4").Value = "NSE:PHOENIXLTD"
$ws.Range("F4").Value = "NSE:LICI"
$ws.Range("B5").Value = "NSE:ADVANIHOTR"
$ws.Range("C5").Value = "NSE:LAL"
$ws.Range("F5").Value = "NSE:MCX"
$ws.Range("B6").Value = "NSE:AMBER"
$ws.Range("C6").Value = "NSE:MTNL"
$ws.Range("F6").Value = "NSE:NMDC"
$ws.Range("B7").Value = "NSE:ARIHANTCAP"
$ws.Range("B8").Value = "NSE:ARMANFIN"
$ws.Range("B9").Value = "NSE:ASHIANA"
$ws.Range("B10").Value = "NSE:ATALREAL"
$ws.Range("B11").Value = "NSE:ATL"
$ws.Range("B12").Value = "NSE:BCLIND"
$ws.Range("B13").Value = "NSE:BSLSENETFG"
$ws.Range("B14").Value = "NSE:COALINDIA"
$ws.Range("B15").Value = "NSE:CONSUMBEES"
$ws.Range("B16").Value = "NSE:CSLFINANCE"
$ws.Range("B17").Value = "NSE:DALBHARAT"
$ws.Range("B18").Value = "NSE:DATAMATICS"
$ws.Range("B19").Value = "NSE:DIVOPPBEES"
$ws.Range("B20").Value = "NSE:DLINKINDIA"
$ws.Range("B21").Value = "NSE:EIHAHOTELS"
$ws.Range("B22").Value = "NSE:ESAFSFB"
$ws.Range("B23").Value = "NSE:FIVESTAR"
$ws.Range("B24").Value = "NSE:FOCUS"
$ws.Range("B25").Value = "NSE:FUSION"
$ws.Range("B26").Value = "NSE:GANESHHOUC"
$ws.Range("B27").Value = "NSE:GATECHDVR"
$ws.Range("B28").Value = "NSE:GEOJITFSL"
$ws.Range("B29").Value = "NSE:HDFCSILVER"
$ws.Range("B30").Value = "NSE:HGS"
$ws.Range("B31").Value = "NSE:INFOMEDIA"
$ws.Range("B32").Value = "NSE:INFRABEES"
$ws.Range("B33").Value = "NSE:ITBEES"
$ws.Range("B34").Value = "NSE:JAMNAAUTO"
$ws.Range("B35").Value = "NSE:JPPOWER"
$ws.Range("B36").Value = "NSE:JSL"
$ws.Range("B37").Value = "NSE:KREBSBIO"
$ws.Range("B38").Value = "NSE:LICNETFN50"
$ws.Range("B39").Value = "NSE:LOWVOL"
$ws.Range("B40").Value = "NSE:MANAKALUCO"
$ws.Range("B41").Value = "NSE:MCX"
$ws.Range("B42").Value = "NSE:MOLDTECH"
$ws.Range("B43").Value = "NSE:MOM100"
$ws.Range("B44").Value = "NSE:MONIFTY500"
$ws.Range("B45").Value = "NSE:NIF100BEES"
$ws.Range("B46").Value = "NSE:NIFTYQLITY"
$ws.Range("B47").Value = "NSE:NMDC"
$ws.Range("B48").Value = "NSE:NPBET"
$ws.Range("B49").Value = "NSE:NRAIL"
$ws.Range("B50").Value = "NSE:NSIL"

# --- Clear cells that no longer have values ---
$ws.Range("C7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()

# --- Append new rows 51-52, copying row formatting from row 9 (plain data row) ---
$ws.Range("A9:F9").Copy($ws.Range("A51:F51"))
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "NSE:ORIENTBELL"

$ws.Range("A9:F9").Copy($ws.Range("A52:F52"))
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "NSE:PAISALO"
